$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ReqPow_AC) - recalculated values
$ws.Range("B2").Value = 88.4545441393966
$ws.Range("C2").Value = 176.9090882787932
$ws.Range("D2").Value = 1101.96296801947
$ws.Range("E2").Value = 1769.090882787932
$ws.Range("F2").Value = 1769.090882787932
$ws.Range("G2").Value = 1590.364055129905
$ws.Range("H2").Value = 1590.364055129905
$ws.Range("I2").Value = 1484.992756861953
$ws.Range("J2").Value = 1484.992756861953
$ws.Range("K2").Value = 176.9090882787932
$ws.Range("L2").Value = 1769.090882787932
$ws.Range("M2").Value = 1769.090882787932
$ws.Range("N2").Value = 1484.992756861953
$ws.Range("O2").Value = 1484.992756861953
$ws.Range("P2").Value = 1484.992756861953
$ws.Range("Q2").Value = 1484.992756861953
$ws.Range("R2").Value = 176.9090882787932
$ws.Range("S2").Value = 176.9090882787932
$ws.Range("T2").Value = 176.9090882787932
$ws.Range("U2").Value = 88.4545441393966

# Row 3 (ReqPow_FC) - recalculated values
$ws.Range("B3").Value = 88.4545441393966
$ws.Range("C3").Value = 176.9090882787932
$ws.Range("D3").Value = 886.5192438262889
$ws.Range("E3").Value = 1437.033505445532
$ws.Range("F3").Value = 1437.033505445532
$ws.Range("G3").Value = 1645.013641234165
$ws.Range("H3").Value = 1645.013641234165
$ws.Range("I3").Value = 1484.992756861953
$ws.Range("J3").Value = 1484.992756861953
$ws.Range("K3").Value = 176.9090882787932
$ws.Range("L3").Value = 1437.033505445532
$ws.Range("M3").Value = 1437.033505445532
$ws.Range("N3").Value = 1484.992756861953
$ws.Range("O3").Value = 1484.992756861953
$ws.Range("P3").Value = 1484.992756861953
$ws.Range("Q3").Value = 1484.992756861953
$ws.Range("R3").Value = 176.9090882787932
$ws.Range("S3").Value = 176.9090882787932
$ws.Range("T3").Value = 176.9090882787932
$ws.Range("U3").Value = 88.4545441393966

# Row 4 (ReqPow_Batt) - recalculated values
$ws.Range("D4").Value = 215.443724193181
$ws.Range("E4").Value = 332.0573773423994
$ws.Range("F4").Value = 332.0573773423994
$ws.Range("G4").Value = -54.64958610426001
$ws.Range("H4").Value = -54.64958610426001
$ws.Range("L4").Value = 332.0573773423994
$ws.Range("M4").Value = 332.0573773423994
